$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting old D:K data to E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy formatting from column E (old D) into new column D so the newly
# inserted cells keep the same number formats / styles as their neighbours
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest reporting period data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 904400
$ws.Range("D9").Value = 464400
$ws.Range("D10").Value = 440100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 40800
$ws.Range("D15").Value = 174900
$ws.Range("D17").Value = 775800
$ws.Range("D18").Value = 128600
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 303600
$ws.Range("D22").Value = 93300
$ws.Range("D23").Value = 35300
$ws.Range("D24").Value = 6200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 29200
$ws.Range("D27").Value = 20200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 20200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 20200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 5600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 163200
$ws.Range("D44").Value = 76300
$ws.Range("D45").Value = 11000
$ws.Range("D46").Value = 256200
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 2171000
$ws.Range("D49").Value = 52400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 72900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2552500
$ws.Range("D57").Value = 54900
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 95800
$ws.Range("D60").Value = 150700
$ws.Range("D61").Value = 1529500
$ws.Range("D62").Value = 23600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1710900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -2263700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 841600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 20200
$ws.Range("D83").Value = 174900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 225900
$ws.Range("D91").Value = -319100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -284900
$ws.Range("D96").Value = -58300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 54100
$ws.Range("D101").Value = "NA"
$ws.Range("D102").Value = -4900

# A few prior-period figures were also revised as part of this refresh
$ws.Range("E89").Value = 201700
$ws.Range("E94").Value = -174500
$ws.Range("E101").Value = "NA"
